$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: shift C1/D1/E1 values (C1<-D1, D1<-E1, E1<-C1's old "max")
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2 data
$ws.Range("C2").Value = "g__Treponema_D"
$ws.Range("D2").Value = "g__Treponema_D"
$ws.Range("E2").Value = 1
